$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.078.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.416.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.22%  "
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.354"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.851.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.014.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.417.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "327.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "
$ws.Range("E24").Value = "  +4.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  +5.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0772"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("E32").Value = "  +3.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.401"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.33"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.74%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  +0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "322.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.96%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "145.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.10%  "
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0516"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.578"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.32%  "
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  -1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.938"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.49%  "
